$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Jonas Valanciunas", "C", "Washington Wizards"),
    @("Malcolm Brogdon", "PG,SG", "Washington Wizards"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Dennis Schröder", "PG", "Brooklyn Nets")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
